$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 14)
$ws.Range("B2").Value = 0.738610351305732
$ws.Range("K2").Value = 0.724772494650826
$ws.Range("L2").Value = 0.77191348422519
$ws.Range("N2").Value = 0.673083606323479

# Row 3 (A3 = 15)
$ws.Range("B3").Value = 0.692201376910046
$ws.Range("K3").Value = 0.661253305020267
$ws.Range("L3").Value = 0.694968829774731
$ws.Range("N3").Value = 0.665790404538871

# Row 4 (A4 = 16)
$ws.Range("B4").Value = 0.682230493525959
$ws.Range("K4").Value = 0.606507362007239
$ws.Range("L4").Value = 0.738494994443583
$ws.Range("N4").Value = 0.636399135871967
